# Updated symbol list on Wed Dec 28 18:28:47 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells are stored as plain text in this workbook, even
# though their contents look numeric. Assigning a numeric-looking string to
# .Value would normally get auto-converted to a real number by Excel, so we
# force the cell to Text format first, write the string, then strip the
# formatting change back off (ClearFormats) so no stray number format is
# left behind on the cell - it ends up as a plain text cell again, matching
# the original layout.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2"  "244.20"
Set-TextValue "D3"  "23.95"
Set-TextValue "D4"  "5.258"
Set-TextValue "D5"  "0.05834"
Set-TextValue "D6"  "6.455"
Set-TextValue "D7"  "3.322"
Set-TextValue "D8"  "0.8081"
Set-TextValue "D9"  "0.8890"
Set-TextValue "D10" "0.1379"
Set-TextValue "D11" "0.07099"
Set-TextValue "D12" "0.03079"
Set-TextValue "D13" "0.03030"
Set-TextValue "D14" "0.09318"
Set-TextValue "D15" "3.817"
Set-TextValue "D16" "0.001538"
Set-TextValue "D17" "0.04712"

Set-TextValue "D18" "0.0006037"
$ws.Range("E18").Value = "17OneONEWorstin24h"

Set-TextValue "D19" "0.006143"
Set-TextValue "D20" "0.001258"
Set-TextValue "D21" "0.004078"
Set-TextValue "D22" "0.00008699"
Set-TextValue "D23" "3.550"
Set-TextValue "D24" "2.173"
Set-TextValue "D25" "0.3187"
Set-TextValue "D26" "0.1312"

Set-TextValue "D40" "0.03841"
Set-TextValue "D41" "0.006299"
Set-TextValue "D43" "0.002600"
Set-TextValue "D44" "0.007007"
Set-TextValue "D45" "0.00005330"

Set-TextValue "D47" "0.5398"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

Set-TextValue "D48" "0.005882"
Set-TextValue "D49" "0.00002100"
Set-TextValue "D50" "0.0002000"
